$wb = $excel.ActiveWorkbook

# Locate the last existing date sheet (2025-08-09) to use as a template,
# then copy it to the end of the workbook to create the new day's sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy([System.Reflection.Missing]::Value, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "2025-08-10"

$data = @(
    ,@(1, "いとこのこ", "いぬちく(著者)", "連載２周年記念マンガ")
    ,@(2, "クセ強彼女は床にいざなう", "須河篤志(著者)", "第14話後半")
    ,@(3, "リビルドワールド", "綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)", "第72話①")
    ,@(4, "ぽんドロイド！ はまさん", "はれやまはれぞう(著者)", "第6話")
    ,@(5, "アイドル辞めるけど結婚してくれますか!?", "三吉汐美(著者)", "休載イラスト")
    ,@(6, "三枝さんはメガネ先輩と恋を描く", "セレビィ量産型(著者)", "第21話後半")
    ,@(7, "モテはるねえ蘆屋くん", "栗原和明(原作) 寺井赤音(作画)", "第16話後半")
    ,@(8, "理想の彼女", "もりまりも(著者)", "第23話")
    ,@(9, "田舎の黒ギャルJKと結婚しました", "カヅチ(著者)", "第18話前半")
    ,@(10, "オトナを知りたい女友達", "望公太(原作) ぷよちゃ(作画)", "第8話 前半")
    ,@(11, "王子様の友達", "すけろく(著者)", "第29話")
    ,@(12, "転生してあらゆるモノに好かれながら異世界で好きな事をして生きて行く", "都尾琉(漫画) 御峰。(原作)", "第27話②")
    ,@(13, "姫騎士は蛮族の嫁", "コトバノリアキ", "第53話①`t焔匠は遺産の番人")
    ,@(14, "悪人面したＢ級冒険者 主人公とその幼馴染たちのパパになる", "こげめ(著者) えんじ(原作) ハラカズヒロ(キャラクター原案)", "第17話-1：「違法奴隷商討伐」")
    ,@(15, "ミルク搾りハンターの異世界搾乳記～農家の冴えない男があらゆる種族の地区Bを弄び虜にする～", "空詠大智(著者)", "第17話後半")
    ,@(16, "リアリスト魔王による聖域なき異世界改革", "鈴木マナツ(漫画) 羽田遼亮(原作) ゆーげん(キャラクターデザイン) ひたきゆう(キャラクターデザイン)", "第68幕①")
    ,@(17, "勇者に全部奪われた俺は勇者の母親とパーティを組みました！", "久遠まこと(著者) 石のやっさん(原作)", "第29話")
    ,@(18, "ダークサモナーとデキている", "車王(著者)", "第73話")
    ,@(19, "チュートリアルが始まる前に ボスキャラ達を破滅させない為に俺ができる幾つかの事", "横山コウヂ(漫画) 髙橋炬燵(原作) カカオ・ランタン(キャラクターデザイン)", "第14話②")
    ,@(20, "序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する", "作画：マエD 原作：新人", "第5話(3)")
    ,@(21, "まんきつしたい常連さん", "しんみりん(著者)", "第46話後編")
    ,@(22, "みつばものがたり 呪いの少女と死の輪舞《ロンド》", "堤りん(漫画) 七沢またり(原作) EURA(キャラクター原案)", "第11話：勝利の美酒")
    ,@(23, "最強の少年聖騎士、転生者を狩る", "作画：御塩 原作：宇奈木ユラ", "第6話(3)")
    ,@(24, "愚かな天使は悪魔と踊る", "アズマサワヨシ(著者)", "第101話①")
    ,@(25, "最強勇者パーティーは愛が知りたい", "山田肌襦袢", "第28話「最後はこぶしがあればいい」")
    ,@(26, "聖液鍛冶屋のエロランタ", "しげきっくす(著者)", "1巻発売記念イラスト")
    ,@(27, "淫獄団地", "搾精研究所(原作) 丈山雄為(漫画)", "第49話（後編）")
    ,@(28, "アイツノカノジョ", "肉丸", "第52話")
    ,@(29, "ワンパンマン", "原作/ＯＮＥ 作画/村田雄介", "208撃目")
    ,@(30, "よくわからないけれど異世界に転生していたようです", "内々けやき あし カオミン", "第137話 よくわからないけれど脱出するみたいです（１）")
    ,@(31, "美人女上司滝沢さん", "やんBARU(著者)", "第202.5話")
    ,@(32, "解雇された暗黒兵士(30代)のスローなセカンドライフ", "岡沢六十四 るれくちぇ sage・ジョー", "第71話(後編) ダリエルVS.滾り")
    ,@(33, "陰キャの俺が席替えでS級美少女に囲まれたら秘密の関係が始まった。", "星野 星野(原作) バラマツヒトミ(漫画) 黒兎 ゆう(キャラクター原案)", "第4話")
    ,@(34, "ステータス・オール∞（インフィニティ） ∞使いの最強能力者、異世界を自由気ままに暮らします！", "漫画：三津屋みやこ 原作：八又ナガト", "第64話")
    ,@(35, "治癒魔法の間違った使い方 ~戦場を駆ける回復要員~", "九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)", "第81話その2")
    ,@(36, "王様ランキング200話～", "十日草輔（とおかそうすけ）", "第263話")
    ,@(37, "十歳の最強魔導師", "猫月 天乃聖樹", "第2話（前編）")
    ,@(38, "ラスボスたちの隠し仔　～魔王城に転生した元社畜プログラマーは自由気ままに『魔導言語《マジックコード》』を開発する～", "森清士郎(漫画) 熊乃げん骨(原作)", "第15話①")
    ,@(39, "本物のカノジョにしたくなるまで、私で試していいよ。", "wano(漫画) 有丈ほえる(原作) 緋月ひぐれ(キャラクター原案)", "第9話①")
    ,@(40, "おっさん、転生して天才役者になる", "芽々ノ圭(漫画) ほえ太郎(原作)", "第38話①")
    ,@(41, "塔の管理をしてみよう", "盧恩＆雪笠(Friendly Land)(著者) 早秋(原作) 雨神(キャラクター原案)", "第91話後編")
    ,@(42, "魔物喰らい ランキング最下位の冒険者は魔物の力で最強へ", "赤井ハコ(漫画) 緒二葉(原作) とよた瑣織(キャラクターデザイン)", "第9話-2")
    ,@(43, "最強で最速の無限レベルアップ ～スキル【経験値1000倍】と【レベルフリー】でレベル上限の枷が外れた俺は無双する～", "シオヤマ琴 鳥羽田 航 トモゼロ", "第77話 答え")
    ,@(44, "異世界ではじめる二拠点生活 ～空間魔法で王都と田舎をいったりきたり～", "丸山りん(漫画) 錬金王(原作) あんべよしろう(キャラクター原案)", "コミックス第1巻発売告知")
    ,@(45, "終末ツーリング", "さいとー栄(著者)", "第49話 三沢基地　その５②")
    ,@(46, "ノロマ魔法と呼ばれた魔法使いは重力魔法で無双する　～まだ重力の概念のない世界にて、少年は万有引力の王となる～", "神原絵理華(漫画) 一森一輝(原作)", "第18話④")
    ,@(47, "腹ペコ魔王と捕虜勇者！～魔王が俺の部屋に飯を食いに来るんだが～", "梅原うめ(作画) ちょきんぎょ。(原作)", "第2話後半")
    ,@(48, "彼女にしたい女子一位、の隣で見つけたあまりちゃん", "寝巻ネルゾ(漫画) 裕時悠示(原作) たん旦(キャラクター原案)", "コミックス1巻のお知らせ")
    ,@(49, "人生逆転ダンジョン", "尾張 ニコ どすこい花丸様", "第4-1話`t道化師")
    ,@(50, "センパイ、自宅警備員の雇用はいかがですか？", "漫画：コブラサナギ 原作：二上圭 キャラ原案：日向あずり", "第6話前半")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
}

# Restore the originally active sheet/selection so this edit only adds the
# new sheet without perturbing the workbook's active-tab state.
$wb.Worksheets.Item(1).Activate()
